{"js": "// Update each division-problem cell's text to the new value.\n// Each old string is unique within the document, so an exact,\n// case-sensitive search safely targets the single matching run.\nconst replacements = [\n  [\"50\u00f73=16, 2\", \"89\u00f73=29, 2\"],\n  [\"53\u00f72=26, 1\", \"43\u00f79=4, 7\"],\n  [\"26\u00f74=6, 2\", \"35\u00f78=4, 3\"],\n  [\"56\u00f74=14, 0\", \"49\u00f76=8, 1\"],\n  [\"11\u00f79=1, 2\", \"86\u00f79=9, 5\"],\n  [\"84\u00f72=42, 0\", \"62\u00f78=7, 6\"],\n  [\"65\u00f73=21, 2\", \"21\u00f73=7, 0\"],\n  [\"92\u00f79=10, 2\", \"66\u00f73=22, 0\"],\n  [\"28\u00f72=14, 0\", \"37\u00f75=7, 2\"],\n  [\"91\u00f73=30, 1\", \"47\u00f75=9, 2\"],\n  [\"86\u00f73=28, 2\", \"57\u00f79=6, 3\"],\n  [\"53\u00f79=5, 8\", \"25\u00f79=2, 7\"],\n  [\"55\u00f74=13, 3\", \"95\u00f77=13, 4\"],\n  [\"44\u00f79=4, 8\", \"51\u00f73=17, 0\"],\n  [\"39\u00f79=4, 3\", \"29\u00f79=3, 2\"],\n  [\"87\u00f79=9, 6\", \"56\u00f77=8, 0\"],\n  [\"43\u00f73=14, 1\", \"28\u00f75=5, 3\"],\n  [\"98\u00f76=16, 2\", \"33\u00f74=8, 1\"],\n  [\"29\u00f74=7, 1\", \"80\u00f77=11, 3\"],\n  [\"88\u00f77=12, 4\", \"88\u00f73=29, 1\"],\n  [\"50\u00f77=7, 1\", \"26\u00f78=3, 2\"],\n  [\"44\u00f73=14, 2\", \"84\u00f75=16, 4\"],\n  [\"60\u00f78=7, 4\", \"23\u00f77=3, 2\"],\n  [\"71\u00f79=7, 8\", \"60\u00f77=8, 4\"],\n  [\"40\u00f78=5, 0\", \"77\u00f74=19, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, 'Replace');\n}\n\nawait context.sync();\n", "ps1": "# Update each division-problem cell's text to the new value.\n# Each old string is unique within the document, so Find/Replace\n# with exact match safely targets the single matching run.\n$replacements = @(\n    @(\"50\u00f73=16, 2\", \"89\u00f73=29, 2\"),\n    @(\"53\u00f72=26, 1\", \"43\u00f79=4, 7\"),\n    @(\"26\u00f74=6, 2\", \"35\u00f78=4, 3\"),\n    @(\"56\u00f74=14, 0\", \"49\u00f76=8, 1\"),\n    @(\"11\u00f79=1, 2\", \"86\u00f79=9, 5\"),\n    @(\"84\u00f72=42, 0\", \"62\u00f78=7, 6\"),\n    @(\"65\u00f73=21, 2\", \"21\u00f73=7, 0\"),\n    @(\"92\u00f79=10, 2\", \"66\u00f73=22, 0\"),\n    @(\"28\u00f72=14, 0\", \"37\u00f75=7, 2\"),\n    @(\"91\u00f73=30, 1\", \"47\u00f75=9, 2\"),\n    @(\"86\u00f73=28, 2\", \"57\u00f79=6, 3\"),\n    @(\"53\u00f79=5, 8\", \"25\u00f79=2, 7\"),\n    @(\"55\u00f74=13, 3\", \"95\u00f77=13, 4\"),\n    @(\"44\u00f79=4, 8\", \"51\u00f73=17, 0\"),\n    @(\"39\u00f79=4, 3\", \"29\u00f79=3, 2\"),\n    @(\"87\u00f79=9, 6\", \"56\u00f77=8, 0\"),\n    @(\"43\u00f73=14, 1\", \"28\u00f75=5, 3\"),\n    @(\"98\u00f76=16, 2\", \"33\u00f74=8, 1\"),\n    @(\"29\u00f74=7, 1\", \"80\u00f77=11, 3\"),\n    @(\"88\u00f77=12, 4\", \"88\u00f73=29, 1\"),\n    @(\"50\u00f77=7, 1\", \"26\u00f78=3, 2\"),\n    @(\"44\u00f73=14, 2\", \"84\u00f75=16, 4\"),\n    @(\"60\u00f78=7, 4\", \"23\u00f77=3, 2\"),\n    @(\"71\u00f79=7, 8\", \"60\u00f77=8, 4\"),\n    @(\"40\u00f78=5, 0\", \"77\u00f74=19, 1\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
